$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 64
$ws.Cells.Item(4, 6).Value = 397
$ws.Cells.Item(5, 6).Value = 1712
$ws.Cells.Item(6, 6).Value = 725
$ws.Cells.Item(7, 6).Value = 2751
$ws.Cells.Item(8, 6).Value = 2125
$ws.Cells.Item(9, 6).Value = 876
$ws.Cells.Item(10, 6).Value = 2357
$ws.Cells.Item(11, 6).Value = 740
$ws.Cells.Item(12, 6).Value = 6832
$ws.Cells.Item(13, 6).Value = 140
$ws.Cells.Item(16, 6).Value = 1559
$ws.Cells.Item(17, 6).Value = 1356
$ws.Cells.Item(18, 6).Value = 1228
$ws.Cells.Item(19, 6).Value = 107
$ws.Cells.Item(20, 6).Value = 2787
$ws.Cells.Item(21, 6).Value = 2472
$ws.Cells.Item(22, 6).Value = 2472
$ws.Cells.Item(23, 6).Value = 813
$ws.Cells.Item(24, 6).Value = 1133
$ws.Cells.Item(25, 6).Value = 272
$ws.Cells.Item(26, 6).Value = 5472
$ws.Cells.Item(27, 6).Value = 302
$ws.Cells.Item(30, 6).Value = 3814
$ws.Cells.Item(31, 6).Value = 179
$ws.Cells.Item(32, 6).Value = 646
$ws.Cells.Item(33, 6).Value = 1725
$ws.Cells.Item(34, 6).Value = 1086
$ws.Cells.Item(35, 6).Value = 191
$ws.Cells.Item(37, 6).Value = 87
$ws.Cells.Item(38, 6).Value = 289
$ws.Cells.Item(39, 6).Value = 1067
$ws.Cells.Item(40, 6).Value = 426
$ws.Cells.Item(42, 6).Value = 55
$ws.Cells.Item(43, 6).Value = 45
$ws.Cells.Item(44, 6).Value = 113
$ws.Cells.Item(45, 6).Value = 927
$ws.Cells.Item(46, 6).Value = 522
$ws.Cells.Item(47, 6).Value = 50
$ws.Cells.Item(48, 6).Value = 14
$ws.Cells.Item(49, 6).Value = 64
$ws.Cells.Item(50, 6).Value = 96

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(10, 6).Value = 407
$ws.Cells.Item(13, 6).Value = 110
$ws.Cells.Item(14, 6).Value = 967
$ws.Cells.Item(17, 6).Value = 16
$ws.Cells.Item(20, 6).Value = 616
$ws.Cells.Item(21, 6).Value = 278
$ws.Cells.Item(22, 6).Value = 363
$ws.Cells.Item(25, 6).Value = 88
$ws.Cells.Item(28, 6).Value = 316
$ws.Cells.Item(29, 6).Value = 84
$ws.Cells.Item(35, 6).Value = 114
$ws.Cells.Item(37, 6).Value = 223
$ws.Cells.Item(44, 6).Value = 1

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 3305
$ws.Cells.Item(5, 6).Value = 414
$ws.Cells.Item(6, 6).Value = 16
$ws.Cells.Item(7, 6).Value = 1486
$ws.Cells.Item(8, 6).Value = 771
$ws.Cells.Item(9, 6).Value = 419
$ws.Cells.Item(10, 6).Value = 2877
$ws.Cells.Item(11, 6).Value = 340
$ws.Cells.Item(12, 6).Value = 634
$ws.Cells.Item(13, 6).Value = 754
$ws.Cells.Item(14, 6).Value = 1267

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 414
$ws.Cells.Item(3, 6).Value = 1486
$ws.Cells.Item(5, 6).Value = 397
$ws.Cells.Item(6, 6).Value = 419
$ws.Cells.Item(7, 6).Value = 2877
$ws.Cells.Item(8, 6).Value = 1712
$ws.Cells.Item(9, 6).Value = 725
$ws.Cells.Item(10, 6).Value = 2751
$ws.Cells.Item(11, 6).Value = 340
$ws.Cells.Item(12, 6).Value = 2125
$ws.Cells.Item(13, 6).Value = 876
$ws.Cells.Item(14, 6).Value = 2357
$ws.Cells.Item(15, 6).Value = 740
$ws.Cells.Item(16, 6).Value = 6832
$ws.Cells.Item(17, 6).Value = 140
$ws.Cells.Item(18, 6).Value = 634
$ws.Cells.Item(19, 6).Value = 754
$ws.Cells.Item(20, 6).Value = 1559
$ws.Cells.Item(21, 6).Value = 1356
$ws.Cells.Item(22, 6).Value = 1228
$ws.Cells.Item(23, 6).Value = 107
$ws.Cells.Item(24, 6).Value = 1267
$ws.Cells.Item(25, 6).Value = 2787
$ws.Cells.Item(26, 6).Value = 2473
$ws.Cells.Item(27, 6).Value = 363
$ws.Cells.Item(28, 6).Value = 813
$ws.Cells.Item(29, 6).Value = 1133
$ws.Cells.Item(30, 6).Value = 272
$ws.Cells.Item(31, 6).Value = 5472
$ws.Cells.Item(32, 6).Value = 302
$ws.Cells.Item(33, 6).Value = 3814
$ws.Cells.Item(34, 6).Value = 646
$ws.Cells.Item(35, 6).Value = 316
$ws.Cells.Item(36, 6).Value = 1725
$ws.Cells.Item(37, 6).Value = 1086
$ws.Cells.Item(38, 6).Value = 84
$ws.Cells.Item(39, 6).Value = 87
$ws.Cells.Item(40, 6).Value = 289
$ws.Cells.Item(41, 6).Value = 1067
$ws.Cells.Item(42, 6).Value = 426
$ws.Cells.Item(44, 6).Value = 55
$ws.Cells.Item(45, 6).Value = 113
$ws.Cells.Item(46, 6).Value = 927
$ws.Cells.Item(47, 6).Value = 522
$ws.Cells.Item(49, 6).Value = 223
$ws.Cells.Item(50, 6).Value = 223
$ws.Cells.Item(51, 6).Value = 96
